$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the contact name and email (row 4) - "Andy amaya" -> "Steve Head"
$ws.Range("A4").Value = "Steve Head"
$ws.Range("B4").Value = "stevehead@nxglabs.in"

# Update the phone number in C4
$ws.Range("C4").Value = 336746546

# Add a mailto hyperlink on B4 pointing at the new e-mail address, preserving
# the cell's existing (non-hyperlink) formatting/style.
$origStyle = $ws.Range("B4").Style
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:stevehead@nxglabs.in")
$ws.Range("B4").Style = $origStyle

# Move/update the active selection to C6
$ws.Range("C6").Select()
